# Government Revenue Accounting.xlsx
#
# Add a new "GRA-vehbatsubsidy" worksheet to separately track the battery
# cost share of vehicles (#267). The new sheet starts out as a duplicate
# of the existing "GRA-evsubsidy" sheet (same revenue-accounting weight
# structure, same source row reference on "Set Values Here"), renamed and
# placed immediately after "GRA-evsubsidy" in the tab order.

$wb = $excel.ActiveWorkbook

$evSubsidySheet = $wb.Worksheets.Item("GRA-evsubsidy")

# Worksheet.Copy(Before, After) -- passing the source sheet itself as
# "After" drops the copy right after "GRA-evsubsidy", ahead of
# "GRA-elecgensubsidy".
$evSubsidySheet.Copy($null, $evSubsidySheet)

$newSheet = $wb.Worksheets.Item("GRA-evsubsidy (2)")
$newSheet.Name = "GRA-vehbatsubsidy"
